# Slide 1: replace the "fill-in-the-blank" attendance password textbox with a
# green "Attendance password is written on the board" callout box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the old "TextBox 1" shape (id=2): "Today's Attendance password / ___________"
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "TextBox 1") {
        $s.Shapes.Item($i).Delete()
        break
    }
}

# Add the new green rounded-outline rectangle (appended after "Picture 2" in z-order)
$shp = $s.Shapes.AddShape(1, 436.15385826771654, 16.510708661417322, 263.4984251968504, 58.16251968503937)
$shp.Name = "Shape 54"

$shp.Fill.ForeColor.RGB = 0x59BB9B
$shp.Line.ForeColor.RGB = 0x418871
$shp.Line.Weight = 2
$shp.Line.CapStyle = 1

$shp.TextFrame.WordWrap = -1
$shp.TextFrame.MarginLeft = 0
$shp.TextFrame.MarginTop = 0
$shp.TextFrame.MarginRight = 0
$shp.TextFrame.MarginBottom = 0
$shp.TextFrame.AutoSize = 1

$tr = $shp.TextFrame.TextRange
$tr.Text = " Attendance password    "
$tr.Font.Size = 24
$tr.Font.Italic = -1
$tr.Font.Color.RGB = 0xFFFFFF

$r2 = $tr.InsertAfter("`r is written on the board")
$r2.Font.Size = 24
$r2.Font.Italic = -1
$r2.Font.Color.RGB = 0xFFFFFF

# Re-assert the exact autofit height (engine recomputes it on text assignment).
$shp.Height = 58.16251968503937
